$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New text for rows 20, 21 and 22 (replacing the old
# GenFigSourcesOfVarPairWise.R / RandomForestFunc.R / RandomForestGroupVarImp.R
# entries with three newly verified functions).
$ws.Range("A20").Value = "GenFigSourcesOfVarVP.R"
$ws.Range("B20").Value = "In general, there is not much logic. There is a risk of mixing samples up - it is however unlikely that an unmatchning design matrix and data matrix are used together - they must have the same number of samples. For selection of LM22S, the results look as expected, with a higher explained variance for cell subType. For variance explained vs gene expression, the residuals look as expected, so it most likely works. Outlier check: The code is fairly straight forward, and the results look somewhat as expected."
$ws.Rows.Item(20).RowHeight = 75

$ws.Range("A21").Value = "GenDeconvData.R"
$ws.Range("B21").Value = "The code has not been formally verified, but it is fairly simple and generates the expected output. The end results of the deconvolution also very much looks as expected, with the internal lab 4 having the least relative error, and so forth."
$ws.Rows.Item(21).RowHeight = 45

$ws.Range("A22").Value = "GenFigDeconv.R"
$ws.Range("B22").Value = "The code has not been formally verified, but it is fairly simple and generates the expected plots. The results of the deconvolution also very much looks as expected, with the internal lab 4 having the least relative error, and so forth."
$ws.Rows.Item(22).RowHeight = 45

# Update the current selection/view to reflect where the user ended up working.
$ws.Range("B23").Select()
$excel.ActiveWindow.ScrollRow = 13
